# Time Log.xlsx — fill in the previously-blank row 100 on Sheet1 with a
# real "Coding" time entry, then move the viewport/selection down to
# reflect the newly added row (matches the author's edit captured in the
# commit "Removed unneeded drawables. Switched to other Iconify icons.
# Fixed issue with DeviceListActivity.").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 100 was previously just a blank placeholder row (only A100/E100 had
# styles, with E100 evaluating to "" via the shared formula). Populate it
# with a genuine entry: date 2014-10-27, 2:00 PM - 2:25 PM, no
# interruption, Activity = Coding.
$ws.Range("A100").Value = 41939
$ws.Range("B100").Value = 0.58333333333333337
$ws.Range("C100").Value = 0.60069444444444442
$ws.Range("D100").Value = 0
$ws.Range("F100").Value = "Coding"

# E100 already carries the shared formula (si="0") inherited from E4;
# touching neighboring cells is enough to make it recompute, but force a
# full recalculation so every dependent (E100, E104 total, Sheet2 SUMIFs,
# percentages) is refreshed before save.
$excel.CalculateFullRebuild()

# Reflect the edit in the UI state: the author scrolled further down the
# sheet and left the selection on C101 (just under the new row).
$ws.Activate()
$ws.Range("C101").Select()
